# Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b
#
# 1) Bump the "Date" metadata field on the Metadata sheet.
# 2) Add a new "Mapping: Spécification métier vers l'extension ROR
#    AvailableTimeNumberDaysofWeek" column (AL) on the Elements sheet,
#    with a value ("jourSemaine") only on the Extension.value[x] row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh Date value -------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: new mapping column --------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Carry over the header style (s="1") and the data-row style (s="2")
# from the existing last column (AK) onto the new column (AL) before
# writing any content into it.
$ws.Range("AK1").Copy()
$ws.Range("AL1").PasteSpecial(-4122)

$ws.Range("AK2:AK6").Copy()
$ws.Range("AL2:AL6").PasteSpecial(-4122)

# Column header
$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR AvailableTimeNumberDaysofWeek"

# Only the Extension.value[x] row (row 6) carries a mapping value; the
# other data rows (2-5) stay blank, same as the other mapping columns.
$ws.Range("AL6").Value = "jourSemaine"

# Column width for the new column, to match the other wide text columns.
$ws.Columns.Item(38).ColumnWidth = 89.8
